$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "CasesTab" row's Cypher query (B2) is being trimmed: it no longer
# returns the trailing `Cohort` column, so drop the last RETURN line
# (and the now-trailing comma on the prior line).
$newCasesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
 MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
	WHERE s.clinical_study_designation IN ['UBC01'] and diag.stage_of_disease in [ 'T3N0M1', 'Not Applicable']  OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $newCasesQuery

# The row shrank by one wrapped line once the Cohort clause was removed,
# so its autofit height drops from 304.5 to 290 (matching rows 3 & 4).
$ws.Rows.Item(2).RowHeight = 290

# Selection moved off the old C4:E4 block onto B2 (the edited cell).
$ws.Range("B2").Select() | Out-Null
